$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Patient")

$ws.Range("A2").Value = "28543365-fd53-492b-81dc-543a11595f3a"
$ws.Range("B2").Value = "VITC413495"
$ws.Range("C2").Value = "Katharine"
$ws.Range("E2").Value = "Paggetti"
$ws.Range("F2").Value = "Mattolini"
$ws.Range("G2").Value = 24055
$ws.Range("I2").Value = "other"
